# "Generate Report for Handoff"
#
# The handoff-status report has three sheets (Overview, zh-cn, de-de), each
# backed by a table, summarizing the localization pipeline state for the
# "ecbcfea4-fd0a-4ec3-8a4f-33a4de03a97f.md" source file (row 5 on every
# sheet). A fresh handoff xliff was generated for that file, so its
# timestamp columns need to be refreshed on all three sheets:
#
#   - Overview!G5  ("Latest HO Xliff Generate Date")      -> 2016-11-09 05:36:22
#   - zh-cn!H5     ("Latest Handoff Datetime")             -> 2016-11-09 05:36:09
#   - de-de!H5     ("Latest Handoff Datetime")             -> 2016-11-09 05:36:22
#
# These cells are formatted with a date/time number format but store their
# value as literal text (shared string), so we assign plain strings rather
# than dates.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-11-09 05:36:22"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-11-09 05:36:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-11-09 05:36:22"
